$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-10-21 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-10-22 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("80-76=", $true, $false, $false, $false, $false, $true, 1, $false, "74-43=", 2) | Out-Null
$d.Content.Find.Execute("60+10=", $true, $false, $false, $false, $false, $true, 1, $false, "19+77=", 2) | Out-Null
$d.Content.Find.Execute("40+15=", $true, $false, $false, $false, $false, $true, 1, $false, "64-11=", 2) | Out-Null
$d.Content.Find.Execute("1+41=", $true, $false, $false, $false, $false, $true, 1, $false, "20+19=", 2) | Out-Null
$d.Content.Find.Execute("98-11=", $true, $false, $false, $false, $false, $true, 1, $false, "37+49=", 2) | Out-Null
$d.Content.Find.Execute("70+28=", $true, $false, $false, $false, $false, $true, 1, $false, "45+37=", 2) | Out-Null
$d.Content.Find.Execute("50-24=", $true, $false, $false, $false, $false, $true, 1, $false, "53-1=", 2) | Out-Null
$d.Content.Find.Execute("46-17=", $true, $false, $false, $false, $false, $true, 1, $false, "33+56=", 2) | Out-Null
$d.Content.Find.Execute("7+75=", $true, $false, $false, $false, $false, $true, 1, $false, "5+66=", 2) | Out-Null
$d.Content.Find.Execute("6+76=", $true, $false, $false, $false, $false, $true, 1, $false, "65-24=", 2) | Out-Null
$d.Content.Find.Execute("38-13=", $true, $false, $false, $false, $false, $true, 1, $false, "37+41=", 2) | Out-Null
$d.Content.Find.Execute("0+10=", $true, $false, $false, $false, $false, $true, 1, $false, "32-19=", 2) | Out-Null
$d.Content.Find.Execute("22+20=", $true, $false, $false, $false, $false, $true, 1, $false, "73-43=", 2) | Out-Null
$d.Content.Find.Execute("27-20=", $true, $false, $false, $false, $false, $true, 1, $false, "75+12=", 2) | Out-Null
$d.Content.Find.Execute("59+17=", $true, $false, $false, $false, $false, $true, 1, $false, "55-3=", 2) | Out-Null
$d.Content.Find.Execute("57+37=", $true, $false, $false, $false, $false, $true, 1, $false, "75+4=", 2) | Out-Null
$d.Content.Find.Execute("99-78=", $true, $false, $false, $false, $false, $true, 1, $false, "58-5=", 2) | Out-Null
$d.Content.Find.Execute("49+18=", $true, $false, $false, $false, $false, $true, 1, $false, "7+58=", 2) | Out-Null
$d.Content.Find.Execute("35+32=", $true, $false, $false, $false, $false, $true, 1, $false, "96-80=", 2) | Out-Null
$d.Content.Find.Execute("33-20=", $true, $false, $false, $false, $false, $true, 1, $false, "91-51=", 2) | Out-Null
$d.Content.Find.Execute("3+2=", $true, $false, $false, $false, $false, $true, 1, $false, "64-20=", 2) | Out-Null
$d.Content.Find.Execute("12-5=", $true, $false, $false, $false, $false, $true, 1, $false, "93-66=", 2) | Out-Null
$d.Content.Find.Execute("70+22=", $true, $false, $false, $false, $false, $true, 1, $false, "98-63=", 2) | Out-Null
$d.Content.Find.Execute("56+19=", $true, $false, $false, $false, $false, $true, 1, $false, "99-50=", 2) | Out-Null
$d.Content.Find.Execute("13+29=", $true, $false, $false, $false, $false, $true, 1, $false, "0+56=", 2) | Out-Null
$d.Content.Find.Execute("15+79=", $true, $false, $false, $false, $false, $true, 1, $false, "80-16=", 2) | Out-Null
$d.Content.Find.Execute("62-14=", $true, $false, $false, $false, $false, $true, 1, $false, "43+56=", 2) | Out-Null
$d.Content.Find.Execute("64-45=", $true, $false, $false, $false, $false, $true, 1, $false, "28-22=", 2) | Out-Null
$d.Content.Find.Execute("33+41=", $true, $false, $false, $false, $false, $true, 1, $false, "61+26=", 2) | Out-Null
$d.Content.Find.Execute("46+44=", $true, $false, $false, $false, $false, $true, 1, $false, "84+15=", 2) | Out-Null
$d.Content.Find.Execute("76-60=", $true, $false, $false, $false, $false, $true, 1, $false, "72-1=", 2) | Out-Null
$d.Content.Find.Execute("36+36=", $true, $false, $false, $false, $false, $true, 1, $false, "26+55=", 2) | Out-Null
$d.Content.Find.Execute("31-10=", $true, $false, $false, $false, $false, $true, 1, $false, "41-40=", 2) | Out-Null
$d.Content.Find.Execute("4+6=", $true, $false, $false, $false, $false, $true, 1, $false, "25+10=", 2) | Out-Null
$d.Content.Find.Execute("39+31=", $true, $false, $false, $false, $false, $true, 1, $false, "3-0=", 2) | Out-Null
$d.Content.Find.Execute("76-38=", $true, $false, $false, $false, $false, $true, 1, $false, "83-0=", 2) | Out-Null
$d.Content.Find.Execute("57+12=", $true, $false, $false, $false, $false, $true, 1, $false, "3+29=", 2) | Out-Null
$d.Content.Find.Execute("66-38=", $true, $false, $false, $false, $false, $true, 1, $false, "28+44=", 2) | Out-Null
$d.Content.Find.Execute("58-35=", $true, $false, $false, $false, $false, $true, 1, $false, "50+2=", 2) | Out-Null
$d.Content.Find.Execute("88-35=", $true, $false, $false, $false, $false, $true, 1, $false, "40+8=", 2) | Out-Null
$d.Content.Find.Execute("62+30=", $true, $false, $false, $false, $false, $true, 1, $false, "84-32=", 2) | Out-Null
$d.Content.Find.Execute("57-53=", $true, $false, $false, $false, $false, $true, 1, $false, "96-59=", 2) | Out-Null
$d.Content.Find.Execute("93-76=", $true, $false, $false, $false, $false, $true, 1, $false, "11+67=", 2) | Out-Null
$d.Content.Find.Execute("86-3=", $true, $false, $false, $false, $false, $true, 1, $false, "74-17=", 2) | Out-Null
$d.Content.Find.Execute("1+30=", $true, $false, $false, $false, $false, $true, 1, $false, "21-3=", 2) | Out-Null
$d.Content.Find.Execute("39+36=", $true, $false, $false, $false, $false, $true, 1, $false, "10+50=", 2) | Out-Null
$d.Content.Find.Execute("40-24=", $true, $false, $false, $false, $false, $true, 1, $false, "6+93=", 2) | Out-Null
$d.Content.Find.Execute("6+62=", $true, $false, $false, $false, $false, $true, 1, $false, "7+45=", 2) | Out-Null
$d.Content.Find.Execute("0+34=", $true, $false, $false, $false, $false, $true, 1, $false, "76-46=", 2) | Out-Null
$d.Content.Find.Execute("14+16=", $true, $false, $false, $false, $false, $true, 1, $false, "31+47=", 2) | Out-Null
$d.Content.Find.Execute("72-28=", $true, $false, $false, $false, $false, $true, 1, $false, "92-89=", 2) | Out-Null
$d.Content.Find.Execute("71-4=", $true, $false, $false, $false, $false, $true, 1, $false, "19+55=", 2) | Out-Null
$d.Content.Find.Execute("19+35=", $true, $false, $false, $false, $false, $true, 1, $false, "53-42=", 2) | Out-Null
$d.Content.Find.Execute("41+52=", $true, $false, $false, $false, $false, $true, 1, $false, "49-17=", 2) | Out-Null
$d.Content.Find.Execute("77-76=", $true, $false, $false, $false, $false, $true, 1, $false, "77-47=", 2) | Out-Null
$d.Content.Find.Execute("51-29=", $true, $false, $false, $false, $false, $true, 1, $false, "14+8=", 2) | Out-Null
$d.Content.Find.Execute("49+4=", $true, $false, $false, $false, $false, $true, 1, $false, "32+20=", 2) | Out-Null
$d.Content.Find.Execute("66-48=", $true, $false, $false, $false, $false, $true, 1, $false, "29-1=", 2) | Out-Null
$d.Content.Find.Execute("56-1=", $true, $false, $false, $false, $false, $true, 1, $false, "59-3=", 2) | Out-Null
$d.Content.Find.Execute("80-36=", $true, $false, $false, $false, $false, $true, 1, $false, "16+5=", 2) | Out-Null
$d.Content.Find.Execute("19+18=", $true, $false, $false, $false, $false, $true, 1, $false, "21+61=", 2) | Out-Null
$d.Content.Find.Execute("90-8=", $true, $false, $false, $false, $false, $true, 1, $false, "46-39=", 2) | Out-Null
$d.Content.Find.Execute("80-21=", $true, $false, $false, $false, $false, $true, 1, $false, "22-20=", 2) | Out-Null
$d.Content.Find.Execute("69-5=", $true, $false, $false, $false, $false, $true, 1, $false, "17+11=", 2) | Out-Null
$d.Content.Find.Execute("48+45=", $true, $false, $false, $false, $false, $true, 1, $false, "20+33=", 2) | Out-Null
$d.Content.Find.Execute("72-22=", $true, $false, $false, $false, $false, $true, 1, $false, "89-88=", 2) | Out-Null
$d.Content.Find.Execute("67-25=", $true, $false, $false, $false, $false, $true, 1, $false, "1+48=", 2) | Out-Null
$d.Content.Find.Execute("73+26=", $true, $false, $false, $false, $false, $true, 1, $false, "27+55=", 2) | Out-Null
$d.Content.Find.Execute("37-10=", $true, $false, $false, $false, $false, $true, 1, $false, "59-27=", 2) | Out-Null
$d.Content.Find.Execute("40-14=", $true, $false, $false, $false, $false, $true, 1, $false, "81-27=", 2) | Out-Null
$d.Content.Find.Execute("83-1=", $true, $false, $false, $false, $false, $true, 1, $false, "91-76=", 2) | Out-Null
$d.Content.Find.Execute("51-5=", $true, $false, $false, $false, $false, $true, 1, $false, "52+36=", 2) | Out-Null
$d.Content.Find.Execute("90-79=", $true, $false, $false, $false, $false, $true, 1, $false, "66+25=", 2) | Out-Null
$d.Content.Find.Execute("91-16=", $true, $false, $false, $false, $false, $true, 1, $false, "24+73=", 2) | Out-Null
$d.Content.Find.Execute("54+21=", $true, $false, $false, $false, $false, $true, 1, $false, "99-53=", 2) | Out-Null
$d.Content.Find.Execute("38+59=", $true, $false, $false, $false, $false, $true, 1, $false, "11+50=", 2) | Out-Null
$d.Content.Find.Execute("35+0=", $true, $false, $false, $false, $false, $true, 1, $false, "35+63=", 2) | Out-Null
$d.Content.Find.Execute("7+22=", $true, $false, $false, $false, $false, $true, 1, $false, "91-36=", 2) | Out-Null
$d.Content.Find.Execute("80-17=", $true, $false, $false, $false, $false, $true, 1, $false, "64-2=", 2) | Out-Null
$d.Content.Find.Execute("39+28=", $true, $false, $false, $false, $false, $true, 1, $false, "77+4=", 2) | Out-Null
$d.Content.Find.Execute("86-26=", $true, $false, $false, $false, $false, $true, 1, $false, "86-48=", 2) | Out-Null
$d.Content.Find.Execute("28+66=", $true, $false, $false, $false, $false, $true, 1, $false, "80-35=", 2) | Out-Null
$d.Content.Find.Execute("88-20=", $true, $false, $false, $false, $false, $true, 1, $false, "76-5=", 2) | Out-Null
$d.Content.Find.Execute("35+39=", $true, $false, $false, $false, $false, $true, 1, $false, "9+34=", 2) | Out-Null
$d.Content.Find.Execute("98-69=", $true, $false, $false, $false, $false, $true, 1, $false, "26+30=", 2) | Out-Null
$d.Content.Find.Execute("55-53=", $true, $false, $false, $false, $false, $true, 1, $false, "35+21=", 2) | Out-Null
$d.Content.Find.Execute("25+16=", $true, $false, $false, $false, $false, $true, 1, $false, "38-11=", 2) | Out-Null
$d.Content.Find.Execute("21+59=", $true, $false, $false, $false, $false, $true, 1, $false, "35-25=", 2) | Out-Null
$d.Content.Find.Execute("70-50=", $true, $false, $false, $false, $false, $true, 1, $false, "31-25=", 2) | Out-Null
$d.Content.Find.Execute("72-61=", $true, $false, $false, $false, $false, $true, 1, $false, "23+44=", 2) | Out-Null
$d.Content.Find.Execute("8+27=", $true, $false, $false, $false, $false, $true, 1, $false, "86-74=", 2) | Out-Null
$d.Content.Find.Execute("85-68=", $true, $false, $false, $false, $false, $true, 1, $false, "58+35=", 2) | Out-Null
$d.Content.Find.Execute("47-39=", $true, $false, $false, $false, $false, $true, 1, $false, "60-7=", 2) | Out-Null
$d.Content.Find.Execute("48-17=", $true, $false, $false, $false, $false, $true, 1, $false, "37-31=", 2) | Out-Null
$d.Content.Find.Execute("58-17=", $true, $false, $false, $false, $false, $true, 1, $false, "49+40=", 2) | Out-Null
$d.Content.Find.Execute("71-35=", $true, $false, $false, $false, $false, $true, 1, $false, "46-1=", 2) | Out-Null
$d.Content.Find.Execute("27+10=", $true, $false, $false, $false, $false, $true, 1, $false, "31+56=", 2) | Out-Null
$d.Content.Find.Execute("85+8=", $true, $false, $false, $false, $false, $true, 1, $false, "59+7=", 2) | Out-Null
$d.Content.Find.Execute("7-6=", $true, $false, $false, $false, $false, $true, 1, $false, "96-95=", 2) | Out-Null
$d.Content.Find.Execute("84-66=", $true, $false, $false, $false, $false, $true, 1, $false, "95-12=", 2) | Out-Null
